$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - AGE_BASE
$ws.Range("C4").Value = "Age at exposure measure [years]"

# Row 8 - TOT_PA_QX
$ws.Range("C8").Value = "Physical activity from questionnaire data [MET-hr/day]"
$ws.Range("I8").Value = "MET-min/week"

# Row 9 - TOT_PA_AC
$ws.Range("C9").Value = "Physical activity from accelerometry data [MET-hr/day]"
$ws.Range("I9").ClearContents()
$ws.Range("K9").Value = "unavailable"

# Row 10 - SMOKE_ST
$ws.Range("F10").Value = "impossible"
$ws.Range("G10").Value = "impossible"
$ws.Range("H10").Value = "impossible"
$ws.Range("I10").Value = "unavailable"
$ws.Range("J10").Value = "impossible"
$ws.Range("K10").Value = "unavailable"

# Row 11 - TOBACCO_PY
$ws.Range("C11").Value = "Cumulative lifetime tobacco exposure  [pack years]"
$ws.Range("F11").Value = "impossible"
$ws.Range("G11").Value = "impossible"
$ws.Range("H11").Value = "impossible"
$ws.Range("I11").ClearContents()
$ws.Range("J11").Value = "impossible"
$ws.Range("K11").Value = "unavailable"

# Row 12 - TOBACCO_D
$ws.Range("C12").Value = "Amount of daily tobacco smoked [g/day]"
$ws.Range("D12").Value = "decimal"
$ws.Range("F12").Value = "impossible"
$ws.Range("G12").Value = "impossible"
$ws.Range("H12").Value = "impossible"
$ws.Range("I12").ClearContents()
$ws.Range("J12").Value = "impossible"
$ws.Range("K12").Value = "unavailable"

# Row 13 - AGE_SMOKE_QUIT
$ws.Range("C13").Value = "Age at time of quitting smoking [years]"
$ws.Range("D13").Value = "decimal"
$ws.Range("F13").Value = "impossible"
$ws.Range("G13").Value = "impossible"
$ws.Range("H13").Value = "impossible"
$ws.Range("I13").ClearContents()
$ws.Range("J13").Value = "impossible"
$ws.Range("K13").Value = "unavailable"

# Row 14 - MED_SUPPL
$ws.Range("I14").ClearContents()
$ws.Range("J14").Value = "impossible"
$ws.Range("K14").Value = "unavailable"

# Row 16 - HRT
$ws.Range("I16").Value = "no HRT for all included women due to inclusion criteria, new variable could be created coded such as HRT=0"

# Row 17 - CONTRACEPTIVE
$ws.Range("C17").Value = "Use of contraceptive pills or injections [years]"
$ws.Range("I17").Value = "all included women did not take contraceptives due to inclusion criteria, new variable could be created coding all participants =0"

# Row 18 - LIVE_BIRTHS
$ws.Range("C18").Value = "Number of live births given [Nr. of birth]"
$ws.Range("I18").ClearContents()
$ws.Range("J18").Value = "impossible"
$ws.Range("K18").Value = "unavailable"

# Row 19 - AGE_FIRST_BIRTH
$ws.Range("C19").Value = "Age at the first given birth [years]"
$ws.Range("D19").Value = "decimal"

# Row 24 - PREV_DIAB
$ws.Range("I24").Value = "due to inclusion criteria, all participants were healthy at time of reruitment without a history of diabetes; new variable could be created coded =0 for all participants"

# Row 27 - PREV_CANCER
$ws.Range("I27").Value = "due to inclusion criteria, all participants were healthy at time of reruitment without a history of cancer; new variable could be created coded =0 for all participants"

# Row 28 - FAM1_CHD_STROKE
$ws.Range("I28").Value = "unavailable"
$ws.Range("J28").Value = "impossible"
$ws.Range("K28").Value = "unavailable"

# Row 35 - MELANOMA_SCREEN
$ws.Range("C35").Value = "Screening, skin cancer"

# Row 36 - MAMMO_SCREEN
$ws.Range("C36").Value = "Screening, mammography"

# Row 37 - CERVICAL_SCREEN
$ws.Range("C37").Value = "Screening cervical, smear test"

# Row 38 - MED_STAT
$ws.Range("I38").Value = "due to inclusion criteria, all participants were healthy at time of reruitment without a history of diabetes; new variable could be created coded =0 for all participants"
$ws.Range("J38").Value = "complete"
$ws.Range("K38").Value = "compatible"

# Row 40 - INC_CVD
$ws.Range("I40").Value = "impossible"
$ws.Range("J40").Value = "impossible"
$ws.Range("K40").Value = "unavailable"

# Row 41 - AGE_CVD
$ws.Range("C41").Value = "Age at diagnosis of CVD [years]"
$ws.Range("D41").Value = "decimal"

# Row 43 - AGE_ANGINA
$ws.Range("C43").Value = "Age at diagnosis of angina pectoris [years]"
$ws.Range("D43").Value = "decimal"

# Row 45 - AGE_MI
$ws.Range("C45").Value = "Age at diagnosis of myocardial infarction [years]"
$ws.Range("D45").Value = "decimal"

# Row 47 - AGE_STR
$ws.Range("C47").Value = "Age at diagnosis of stroke [years]"
$ws.Range("D47").Value = "decimal"

# Row 49 - AGE_ISC_STR
$ws.Range("C49").Value = "Age at diagnosis of cerebral infarction (ischaemic stroke) [years]"
$ws.Range("D49").Value = "decimal"

# Row 51 - AGE_HAEMO_STR
$ws.Range("C51").Value = "Age at diagnosis of haemorrhagic stroke [years]"
$ws.Range("D51").Value = "decimal"

# Row 53 - AGE_HYP
$ws.Range("C53").Value = "Age at diagnosis of essential hypertension [years]"
$ws.Range("D53").Value = "decimal"

# Row 55 - AGE_HF
$ws.Range("C55").Value = "Age at diagnosis of heart failure [years]"
$ws.Range("D55").Value = "decimal"

# Row 57 - AGE_DIAB2
$ws.Range("C57").Value = "Age at diagnosis of diabetes mellitus type 2 [years]"
$ws.Range("D57").Value = "decimal"

# Row 59 - TYPE_CANCER
$ws.Range("C59").Value = "Type of Cancer (ICD 10, 3 digits,e.g. C18.3)"
$ws.Range("D59").Value = "text"

# Row 60 - AGE_CANCER
$ws.Range("C60").Value = "Age at diagnosis of cancer [years]"
$ws.Range("D60").Value = "decimal"

# Row 62 - AGE_DEATH
$ws.Range("C62").Value = "Age at time of death [years]"
$ws.Range("D62").Value = "decimal"

# Row 65 - AGE_FUP
$ws.Range("C65").Value = "Age at end of follow-up [years]"
